$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Ideal" coverage table: header + 4 data rows with a ratio column
# mirroring the existing C2:C5 pattern (share of SUM(B2:B5)).
$ws.Range("A13").Value = "covered by # of centers"

$ws.Range("A14").Value = 1
$ws.Range("B14").Value = 123456.91
$ws.Range("C14").Formula = "=B14/SUM(B2:B5)"

$ws.Range("A15").Value = 2
$ws.Range("B15").Value = 50806.17
$ws.Range("C15").Formula = "=B15/SUM(B2:B5)"

$ws.Range("A16").Value = 3
$ws.Range("B16").Value = 30096.84
$ws.Range("C16").Formula = "=B16/SUM(B2:B5)"

$ws.Range("A17").Value = 4
$ws.Range("B17").Value = 4484.16
$ws.Range("C17").Formula = "=B17/SUM(B2:B5)"

# Move the active selection to match the edited workbook state
$ws.Range("H15").Select()
